$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($colRef, $val) {
    # Force the value to be stored as text (matches source data which uses
    # inline strings even for numeric-looking figures like "231.27" or "1.840.32").
    # A leading apostrophe forces text entry (stripped from the stored value),
    # then resetting the style back to Normal drops the transient quote-prefix
    # style so the cell keeps using the original (unstyled) format.
    $ws.Range($colRef).Formula = "'" + $val
    $ws.Range($colRef).Style = "Normal"
}

Set-TextValue "D2" "35.511.47"
$ws.Range("E2").Value = "  +3.01%  "

Set-TextValue "D3" "1.840.32"
$ws.Range("E3").Value = "  +1.97%  "

$ws.Range("E4").Value = "  +0.29%  "

Set-TextValue "D5" "231.27"
$ws.Range("E5").Value = "  +2.81%  "

Set-TextValue "D6" "0.611"
$ws.Range("E6").Value = "  +1.08%  "

$ws.Range("E7").Value = "  +0.23%  "

Set-TextValue "D8" "43.92"
$ws.Range("E8").Value = "  +12.16%  "

Set-TextValue "D9" "0.311"
$ws.Range("E9").Value = "  +7.90%  "

Set-TextValue "D10" "0.0702"
$ws.Range("E10").Value = "  +4.81%  "

$ws.Range("E11").Value = "  +2.52%  "

Set-TextValue "D12" "2.106.47"
$ws.Range("E12").Value = "  +1.87%  "

Set-TextValue "D13" "1.831.85"
$ws.Range("E13").Value = "  +1.62%  "

$ws.Range("E14").Value = "  +7.00%  "

$ws.Range("E15").Value = "  +1.77%  "

Set-TextValue "D16" "4.74"
$ws.Range("E16").Value = "  +8.20%  "

Set-TextValue "D17" "35.482.66"
$ws.Range("E17").Value = "  +2.88%  "

Set-TextValue "D18" "70.34"
$ws.Range("E18").Value = "  +3.11%  "

Set-TextValue "D19" "0.0₃0801"
$ws.Range("E19").Value = "  +4.12%  "

Set-TextValue "D20" "244.45"
$ws.Range("E20").Value = "  +1.83%  "

Set-TextValue "D21" "12.05"
$ws.Range("E21").Value = "  +8.19%  "

$ws.Range("E22").Value = "  +14.88%  "

$ws.Range("E23").Value = "  +0.29%  "

$ws.Range("E24").Value = "  +1.61%  "

Set-TextValue "D25" "171.61"
$ws.Range("E25").Value = "  +0.46%  "

$ws.Range("E26").Value = "  +3.18%  "

$ws.Range("E27").Value = "  +1.13%  "

$ws.Range("E28").Value = "  -0.64%  "

Set-TextValue "D29" "1.58"
$ws.Range("E29").Value = "  +28.50%  "

$ws.Range("E30").Value = "  +0.31%  "

Set-TextValue "D31" "3.306.30"
$ws.Range("E31").Value = "  +36.08%  "

$ws.Range("E32").Value = "  +7.49%  "

$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D33" "3.94"
$ws.Range("E33").Value = "  +5.00%  "

$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D34" "4.07"
$ws.Range("E34").Value = "  +5.97%  "

$ws.Range("E35").Value = "  +2.30%  "

Set-TextValue "D36" "95.76"
$ws.Range("E36").Value = "  +16.56%  "

$ws.Range("E37").Value = "  +7.69%  "

$ws.Range("E38").Value = "  +6.49%  "

$ws.Range("B39").Value = "InjectiveProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D39" "15.59"
$ws.Range("E39").Value = "  +11.93%  "

$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D40" "1.350.08"
$ws.Range("E40").Value = "  +3.39%  "

Set-TextValue "D41" "2.45"
$ws.Range("E41").Value = "  +6.25%  "

Set-TextValue "D42" "0.0195"
$ws.Range("E42").Value = "  +4.86%  "

$ws.Range("E43").Value = "  +6.57%  "

$ws.Range("E44").Value = "  +4.58%  "

Set-TextValue "D45" "2.46"
$ws.Range("E45").Value = "  +0.95%  "

$ws.Range("E46").Value = "  +0.87%  "

Set-TextValue "D47" "6.29"
$ws.Range("E47").Value = "  +8.98%  "

Set-TextValue "D48" "0.0519"
$ws.Range("E48").Value = "  +1.01%  "

Set-TextValue "D49" "2.008.58"
$ws.Range("E49").Value = "  +1.99%  "

$ws.Range("E50").Value = "  +0.28%  "

Set-TextValue "D51" "103.54"
$ws.Range("E51").Value = "  +1.04%  "
